$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "This task uses Latin hypercube sampling in order to model the given system of equations in ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This task uses Latin hypercube sampling in order to model the given system of equations in ",
    2)
